$wb = $excel.ActiveWorkbook

# --- tc_01 (sheet1): remove row 3 (the stray q34234/dfsdf row) and move the
# selection to the now-empty row underneath it ---
$ws1 = $wb.Worksheets.Item("tc_01")
$ws1.Range("A3:B3").ClearContents()
$ws1.Range("A3:B3").Select()

# --- add the new "testData" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "testData"

# ZOHO.com login test data
$newSheet.Range("A1").Value = "tc_02"

$newSheet.Range("A2").Value = "username"
$newSheet.Range("B2").Value = "password"

$newSheet.Range("A3").Value = 8744954505
$newSheet.Range("B3").Value = 12123

$newSheet.Range("A5").Value = "tc_01"

$newSheet.Range("A6").Value = "username"
$newSheet.Range("B6").Value = "password"

$newSheet.Range("A7").Value = "esdf"
$newSheet.Range("B7").Value = "sdwerwe"

$newSheet.Range("A8").Value = "sdfsdf"
$newSheet.Range("B8").Value = "qwewer"

$newSheet.Range("A9").Value = "fsfsd"
$newSheet.Range("B9").Value = "dsdfsdf"

$newSheet.Range("C6").Value = "email"
$newSheet.Range("D6").Value = "firstname"
$newSheet.Range("E6").Value = "lastaname"

# column A matches the width used on the other sheets
$newSheet.Columns.Item(1).ColumnWidth = 10.1666666667

# select E6 and make this new sheet the active/visible tab
$newSheet.Range("E6").Select()
$newSheet.Activate()
